# Juno: check in to OLPRODLOC.
# Update the Korean "Chai Tea" tagline list to the localized wording.

$d = $word.ActiveDocument

$replacements = @(
    @{ old = "차이 차 : 삶의 향신료"; new = "차이 티: 인생의 향기" },
    @{ old = "차이 차 : 컵에 맛의 세계"; new = "차이 티: 한 컵에 담아내는 온 세상의 맛" },
    @{ old = "차이 차 : 인도의 마법을 발견"; new = "차이 티: 인도식 매직의 발견" },
    @{ old = "차이 차 : 건강과 즐거움의 완벽한 혼합"; new = "차이 티: 건강과 즐거움의 완벽한 조화" },
    @{ old = "차이 차: 차 그 이상, 삶의 방식"; new = "차이 티: 차 그 이상, 삶의 방식" },
    @{ old = "차이 차: 사계절과 이유를 위한 음료"; new = "차이 티: 모든 계절과 이유를 위한 음료" },
    @{ old = "차이 차 : 당신의 감각에 대한 궁극적 인 방종"; new = "차이 티: 감각에 대한 최고의 관용" },
    @{ old = "차이 차 : 일상에서 달콤한 탈출"; new = "차이 티: 달콤한 일탈" },
    @{ old = "차이 차 : 따뜻함을 공유하고 사랑을 공유합니다."; new = "차이 티: 따뜻함과 사랑을 나누는 방식" },
    @{ old = "차이 차 : 특별한 무언가에 자신을 치료"; new = "차이 티: 자신을 위한 특별한 선물" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
